# Regenerate the "K" column (G) values for save_data sheet.
# These values represent strikeouts (K) recomputed from source data,
# replacing the prior Strike# derived numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 0
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    10 = 1
    11 = 3
    12 = 1
    13 = 1
    15 = 0
    16 = 1
    17 = 1
    18 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
